# Applies the "Added data for ReplyAll and Delegate" change to the
# Transmittals_New sheet (sheet1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transmittals_New")

# --- Header row: O1 was "Action-Level3" -> becomes "DelegateTo",
#     and the old "Action-Level3" header moves one column right to P1
#     (keeping the bold/shaded header formatting).
$ws.Range("O1").Value = "DelegateTo"
$ws.Range("P1").Value = "Action-Level3"
$ws.Range("N1").Copy()
$ws.Range("P1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Rows 4 & 5: their "Submission" value in column O moves to column P.
$ws.Range("O4").ClearContents()
$ws.Range("P4").Value = "Submission"

$ws.Range("O5").ClearContents()
$ws.Range("P5").Value = "Submission"

# --- New row 8: Delegate test data.
$ws.Range("A8").Value = "AutoTestAdmin"
$ws.Range("C8").Value = "New Transmittal from Automation"
$ws.Range("D8").Value = "UnTick"
$ws.Range("E8").Value = "Correspondence"
$ws.Range("F8").Value = "Issued for Review"
$ws.Range("L8").Value = "Delegate- Message for New transmittal"
$ws.Range("M8").Value = "Delegate"
$ws.Range("O8").Value = "AutoTestUser"
$ws.Range("P8").Value = "Submission"

# --- New row 9: ReplyAll test data.
$ws.Range("A9").Value = "AutoTestAdmin"
$ws.Range("B9").Value = "AutoTestUser"
$ws.Range("C9").Value = "New Transmittal from Automation"
$ws.Range("D9").Value = "UnTick"
$ws.Range("E9").Value = "Correspondence"
$ws.Range("F9").Value = "Issued for Review"
$ws.Range("L9").Value = "Reply All- Message for New transmittal"
$ws.Range("M9").Value = "ReplyAll"
$ws.Range("P9").Value = "Submission"

# --- Match the saved selection/active cell recorded in the sheet view.
$ws.Range("C18").Select()
